# Ottawa/Welcome.pptx maintenance edit:
#  - refresh the two cached "datetimeFigureOut" field stamps (slide-layout
#    "Title and Content" date placeholder + the Notes Master date placeholder)
#  - fix the repo owner in the "Not sure?" slide's github.com link
#    (FBoucher -> msdevmtl)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Cached date field on the "Title and Content" slide layout (layout #30
#    on the slide master) -- Date Placeholder shape.
# ---------------------------------------------------------------------------
$titleAndContentLayout = $p.SlideMaster.CustomLayouts.Item(30)
for ($i = 1; $i -le $titleAndContentLayout.Shapes.Count; $i++) {
    $shp = $titleAndContentLayout.Shapes.Item($i)
    if ($shp.Name -eq "Date Placeholder 3") {
        $shp.TextFrame.TextRange.Text = "3/28/2019"
    }
}

# ---------------------------------------------------------------------------
# 2) Cached date field on the Notes Master -- Date Placeholder shape.
# ---------------------------------------------------------------------------
$notesMaster = $p.NotesMaster
for ($i = 1; $i -le $notesMaster.Shapes.Count; $i++) {
    $shp = $notesMaster.Shapes.Item($i)
    if ($shp.Name -eq "Date Placeholder 2") {
        $shp.TextFrame.TextRange.Text = "3/28/2019"
    }
}

# ---------------------------------------------------------------------------
# 3) "Not sure?" slide (slide 3) -- fix the GitHub repo owner in the URL
#    text box from FBoucher to msdevmtl.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(3)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 9") {

        # Make sure the "Not sure?" title run is marked the same as the
        # rest of the refreshed text (dirty="0" in the final deck).
        $titlePara = $shp.TextFrame.TextRange.Paragraphs(1)
        $titlePara.Runs(1).Text = $titlePara.Runs(1).Text

        # The URL lives in paragraph 2 as a single run:
        #   "github.com/FBoucher/GlobalAzureBootcamp-2019"
        # Replace just the "FBoucher/" segment (characters 12-20) with
        # "msdevmtl/" so the run splits into
        #   "github.com/" | "msdevmtl/" | "GlobalAzureBootcamp-2019".
        $urlPara = $shp.TextFrame.TextRange.Paragraphs(2)
        $urlRun = $urlPara.Runs(1)
        $ownerSeg = $urlRun.Characters(12, 9)
        $ownerSeg.Text = "msdevmtl/"
    }
}
